# Update "Password Audit" sheet: fix header row formatting (it was mistakenly
# applied to the first data row instead of the header) and refresh the bcrypt
# hashes in column C for every dealer record with newly generated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the bold/centered style from row 2 (data row) to row 1 (header row).
# Row 1 header cells become bold + centered; row 2 data cells revert to default formatting.
$headerRange = $ws.Range("A1:C1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108  # xlCenter

$dataRow2Range = $ws.Range("A2:C2")
$dataRow2Range.ClearFormats()

# Update the bcrypt hashes in column C (rows 2-31) with the newly generated hashes.
$ws.Range("C2").Value = '$2b$10$UE4Fy6vRoQ3ux6G8EZyJBOc5MMUPvm404K8cZ3bIXkoOFCKxQf2Ku'
$ws.Range("C3").Value = '$2b$10$pGeUWPatady1022jQm31J.Nz6X6s3ZRVPec4qpJHHFtR/gOr/cusq'
$ws.Range("C4").Value = '$2b$10$OfyN0BEZ8iV3w0GlVZ2m7.jajxCP34XKlFiQaqYEHY6Qm12dC3UDy'
$ws.Range("C5").Value = '$2b$10$RwLatoVIvWS/UV0ZOo6x0O5tOseBcUqmC3GAhO.R6TAM7ORH4c3xi'
$ws.Range("C6").Value = '$2b$10$PrPFHrODDQbxOUgA8xTQc.VoPD1UX4f06tRHPlzEfbpT5VsgMWrVK'
$ws.Range("C7").Value = '$2b$10$W/M/oeyADIWz7FsyUpJ64.5e6zpqsJbwsFLJhOI6V6E0BOxSB54aS'
$ws.Range("C8").Value = '$2b$10$Bqo/DwW7n4IZDELnKr6t0.Hz7S/csVuEd274oCqh74tgADNsaGQAS'
$ws.Range("C9").Value = '$2b$10$DNUXiJAxD.t6pgrt6anBaODTQI2o6m/q8BaOCyWFISeI72qAASSRe'
$ws.Range("C10").Value = '$2b$10$1sulnbcPoAiPfHeYlfJXWOc8AxJkK8.aRixzHx4Yvj4Ev7n3nSuAm'
$ws.Range("C11").Value = '$2b$10$n9gcYmfxARo1Yft8DnPt1emT.iY/YYG4N5nme9T0kEu9oPv8aVrgq'
$ws.Range("C12").Value = '$2b$10$u6KC6NjfbDRafhqRvDaW6.0lwHSpZvclk6s.iPmFn.rQ2anP6U0ai'
$ws.Range("C13").Value = '$2b$10$/m5HnznXjuujiyM9WYpKTOUpDjmxTa3R8.58oZ6uoleXG/zPKDwQm'
$ws.Range("C14").Value = '$2b$10$8QpOSn.C3yhdnelfmFSxTua7myop8Z2pOhfvA802BNfe9k/3MhbXS'
$ws.Range("C15").Value = '$2b$10$aWTSqe37E55CK3IfT75XgOMxTbLl0vBwXrasdtdDJuEtSoqJ6jykq'
$ws.Range("C16").Value = '$2b$10$xLcAH/xu10vX01YeQcUQ6eJGeHa95D.ARBvInaHQPCEucV90IhIF.'
$ws.Range("C17").Value = '$2b$10$OwbsCzpu4mhDQpM9lMk/XucpH8Dv4E92.WBVRTNb4XPm1PmmbgqXW'
$ws.Range("C18").Value = '$2b$10$8.t3VLCcyyZUZwej/Vi3MOv92Lkgw22FLQXToBC0I76LGPszhxbP2'
$ws.Range("C19").Value = '$2b$10$IDErbMq4m.f/ZG.m5gnopuveoC72Y0vNNK85r5F9ybkfSZMmnIzHi'
$ws.Range("C20").Value = '$2b$10$Nq8oFMsICbvEfBYaGhzgSuZTyq489XqDhs0znYGiKZuWaZ6O9w2li'
$ws.Range("C21").Value = '$2b$10$dvI/8aG7CoxCOFEDWyfDW.LMDRT4isu6BfKn.YviC0fch7/PTWzgS'
$ws.Range("C22").Value = '$2b$10$3yixvJmtqFESFlDT0jxdLewblxPsqEUd9j.wEd9dLLzEAnRXHPDKu'
$ws.Range("C23").Value = '$2b$10$1Mwd19UzdiXBAbg.V1JQLehEAczA/dF17oJ45n3TC8GPqtIxrb3aO'
$ws.Range("C24").Value = '$2b$10$H9WsuWgkDZ98eheQknpqtelFEWLs.SeMNWPJxRQM9iNBnOCOdOeqm'
$ws.Range("C25").Value = '$2b$10$fdn9o6qlKmpS0OjAkabZAeMYld1E3bmd3cNTc89H3Z2peHP0LPs4i'
$ws.Range("C26").Value = '$2b$10$MzKBeVISBkgmAx48y8aZNOkx8jQKWA8WJe9ne7oz3eLNXq3bR034y'
$ws.Range("C27").Value = '$2b$10$sz53cehQkiWZSey76ApdHOu64x3tcU7cHrWm6KNLhi/E46HhE8772'
$ws.Range("C28").Value = '$2b$10$N7fTIFOp3uHNYOjRAb6IyOyUofuC5BcIyXXp4UQd.jotlHv1vefVS'
$ws.Range("C29").Value = '$2b$10$KdfYcGkH0QyuncShHRZJNu3rFkXGciNdTCJJ836J5PzfcPmPRUpxC'
$ws.Range("C30").Value = '$2b$10$PxTDh8jwFDOUFPrP34htkuZu0bcv6JMsTDTqH1mKkosfhLJ5/j5gC'
$ws.Range("C31").Value = '$2b$10$Z8QAWEluUwdNJHLQkllFVu4gVzob.TOhJ1DKP65TdRZpOFEvPupxu'
